$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "A"
$ws2.Range("B1").Value = "B"
$ws2.Range("D1").Value = "D"

# Row 2
$ws2.Range("A2").Value = "Hello"
$ws2.Range("B2").Value = "World"
$ws2.Range("D2").Value = 12
$ws2.Range("D2").NumberFormat = "0.00"

# Row 3
$ws2.Range("A3").Value = "Foo"
$ws2.Range("B3").Value = "Bar"
$ws2.Range("D3").Value = 123
$ws2.Range("D3").NumberFormat = "0.00"

# Row 4
$ws2.Range("A4").Value = "Extra"
$ws2.Range("B4").Value = "Bar"
$ws2.Range("C4").Value = 123
$ws2.Range("C4").NumberFormat = "0.00"
$ws2.Range("D4").Value = 12
$ws2.Range("D4").NumberFormat = "0.00"

# Added last, matching shared-string insertion order in the target file
$ws2.Range("C2").Value = "Missing"

# Leave the selection on C3, as in the target sheet
$ws2.Range("C3").Select()
